$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: update the "Date" value
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B9").Value = "2025-12-03T10:56:11+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: append the FHIR version "|4.0.1" to ValueSet URLs,
#    Reference(...) type lists and the Quantity {SimpleQuantity} type.
# ---------------------------------------------------------------------------
$wsEl = $wb.Worksheets.Item("Elements")

# Binding Value Set column (AA) -- plain (non-wrapped) values
$wsEl.Range("AA6").Value  = "http://hl7.org/fhir/ValueSet/languages|4.0.1"
$wsEl.Range("AA16").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$wsEl.Range("AA47").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$wsEl.Range("AA24").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$wsEl.Range("AA49").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$wsEl.Range("AA25").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"
$wsEl.Range("AA50").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"
$wsEl.Range("AA27").Value = "http://hl7.org/fhir/ValueSet/body-site|4.0.1"
$wsEl.Range("AA28").Value = "http://hl7.org/fhir/ValueSet/observation-methods|4.0.1"
$wsEl.Range("AA37").Value = "http://hl7.org/fhir/ValueSet/referencerange-meaning|4.0.1"
$wsEl.Range("AA38").Value = "http://hl7.org/fhir/ValueSet/referencerange-appliesto|4.0.1"

# Type(s) column (L) -- values end with a trailing newline in the workbook
$wsEl.Range("L12").Value = "Reference(CarePlan|4.0.1|DeviceRequest|4.0.1|ImmunizationRecommendation|4.0.1|MedicationRequest|4.0.1|NutritionOrder|4.0.1|ServiceRequest|4.0.1)`n"
$wsEl.Range("L13").Value = "Reference(MedicationAdministration|4.0.1|MedicationDispense|4.0.1|MedicationStatement|4.0.1|Procedure|4.0.1|Immunization|4.0.1|ImagingStudy|4.0.1)`n"
$wsEl.Range("L18").Value = "Reference(Resource|4.0.1)`n"
$wsEl.Range("L19").Value = "Reference(Encounter|4.0.1)`n"
$wsEl.Range("L22").Value = "Reference(Practitioner|4.0.1|PractitionerRole|4.0.1|Organization|4.0.1|CareTeam|4.0.1|Patient|4.0.1|RelatedPerson|4.0.1)`n"
$wsEl.Range("L29").Value = "Reference(Specimen|4.0.1)`n"
$wsEl.Range("L30").Value = "Reference(Device|4.0.1|DeviceMetric|4.0.1)`n"
$wsEl.Range("L35").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$wsEl.Range("L36").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$wsEl.Range("L41").Value = "Reference(Observation|4.0.1|QuestionnaireResponse|4.0.1|MolecularSequence|4.0.1)`n"

# ---------------------------------------------------------------------------
# 3) Elements sheet: the content edits above cause Excel's "best fit" column
#    widths to be recalculated on re-save (the IG publishing tool reruns an
#    auto-fit pass over the whole sheet). Reproduce the resulting narrower
#    column widths as closely as this runtime's width model allows.
# ---------------------------------------------------------------------------
$wsEl.Columns.Item(1).ColumnWidth  = 12.666666666666666
$wsEl.Columns.Item(2).ColumnWidth  = 37.5
$wsEl.Columns.Item(3).ColumnWidth  = 37.5
$wsEl.Columns.Item(4).ColumnWidth  = 10.0
$wsEl.Columns.Item(5).ColumnWidth  = 33.333333333333336
$wsEl.Columns.Item(6).ColumnWidth  = 5.0
$wsEl.Columns.Item(7).ColumnWidth  = 3.3333333333333335
$wsEl.Columns.Item(8).ColumnWidth  = 3.8333333333333335
$wsEl.Columns.Item(9).ColumnWidth  = 13.166666666666666
$wsEl.Columns.Item(11).ColumnWidth = 11.5
$wsEl.Columns.Item(20).ColumnWidth = 6.833333333333333
$wsEl.Columns.Item(21).ColumnWidth = 7.666666666666667
$wsEl.Columns.Item(22).ColumnWidth = 13.833333333333334
$wsEl.Columns.Item(23).ColumnWidth = 14.333333333333334
$wsEl.Columns.Item(24).ColumnWidth = 15.5
$wsEl.Columns.Item(25).ColumnWidth = 15.333333333333334
$wsEl.Columns.Item(26).ColumnWidth = 53.833333333333336
$wsEl.Columns.Item(27).ColumnWidth = 50.833333333333336
$wsEl.Columns.Item(28).ColumnWidth = 4.666666666666667
$wsEl.Columns.Item(29).ColumnWidth = 18.666666666666668
$wsEl.Columns.Item(30).ColumnWidth = 16.833333333333332
$wsEl.Columns.Item(31).ColumnWidth = 14.0
$wsEl.Columns.Item(32).ColumnWidth = 11.5
$wsEl.Columns.Item(33).ColumnWidth = 34.333333333333336
$wsEl.Columns.Item(34).ColumnWidth = 8.166666666666666
